$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Ensure cells are formatted as text so numeric-looking strings
# (e.g. "0.732", "13.00") are preserved exactly as authored,
# matching the source inline-string cells instead of being
# auto-converted into floating point numbers by Excel.

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "35.297.74"
$ws.Range("E2").NumberFormat = "@"
$ws.Range("E2").Value = "  +0.48%  "
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "1.905.05"
$ws.Range("E3").NumberFormat = "@"
$ws.Range("E3").Value = "  +0.77%  "
$ws.Range("E4").NumberFormat = "@"
$ws.Range("E4").Value = "  -0.49%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "0.732"
$ws.Range("E5").NumberFormat = "@"
$ws.Range("E5").Value = "  +10.83%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "256.42"
$ws.Range("E6").NumberFormat = "@"
$ws.Range("E6").Value = "  +4.61%  "
$ws.Range("E7").NumberFormat = "@"
$ws.Range("E7").Value = "  -0.46%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "40.70"
$ws.Range("E8").NumberFormat = "@"
$ws.Range("E8").Value = "  -0.65%  "
$ws.Range("E9").NumberFormat = "@"
$ws.Range("E9").Value = "  +6.43%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "52.96"
$ws.Range("E10").NumberFormat = "@"
$ws.Range("E10").Value = "  +0.54%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.0757"
$ws.Range("E11").NumberFormat = "@"
$ws.Range("E11").Value = "  +6.31%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.0988"
$ws.Range("E12").NumberFormat = "@"
$ws.Range("E12").Value = "  -0.52%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "2.181.56"
$ws.Range("E13").NumberFormat = "@"
$ws.Range("E13").Value = "  +0.57%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "13.00"
$ws.Range("E14").NumberFormat = "@"
$ws.Range("E14").Value = "  +8.08%  "
$ws.Range("E15").NumberFormat = "@"
$ws.Range("E15").Value = "  +4.79%  "
$ws.Range("E16").NumberFormat = "@"
$ws.Range("E16").Value = "  +3.65%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "1.921.06"
$ws.Range("E17").NumberFormat = "@"
$ws.Range("E17").Value = "  +1.59%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "35.280.95"
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "75.03"
$ws.Range("E19").NumberFormat = "@"
$ws.Range("E19").Value = "  +3.82%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "0.0₃0847"
$ws.Range("E20").NumberFormat = "@"
$ws.Range("E20").Value = "  +4.00%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "244.50"
$ws.Range("E21").NumberFormat = "@"
$ws.Range("E21").Value = "  +2.22%  "
$ws.Range("E22").NumberFormat = "@"
$ws.Range("E22").Value = "  +5.83%  "
$ws.Range("E23").NumberFormat = "@"
$ws.Range("E23").Value = "  +6.46%  "
$ws.Range("E24").NumberFormat = "@"
$ws.Range("E24").Value = "  -0.32%  "
$ws.Range("E25").NumberFormat = "@"
$ws.Range("E25").Value = "  +6.84%  "
$ws.Range("E26").NumberFormat = "@"
$ws.Range("E26").Value = "  +5.06%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "166.24"
$ws.Range("E27").NumberFormat = "@"
$ws.Range("E27").Value = "  -2.27%  "
$ws.Range("E28").NumberFormat = "@"
$ws.Range("E28").Value = "  +3.99%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "18.83"
$ws.Range("E29").NumberFormat = "@"
$ws.Range("E29").Value = "  +3.21%  "
$ws.Range("E30").NumberFormat = "@"
$ws.Range("E30").Value = "  +5.42%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "4.129.12"
$ws.Range("E31").NumberFormat = "@"
$ws.Range("E31").Value = "  +19.47%  "
$ws.Range("E32").NumberFormat = "@"
$ws.Range("E32").Value = "  +6.89%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "1.64"
$ws.Range("E33").NumberFormat = "@"
$ws.Range("E33").Value = "  +23.96%  "
$ws.Range("E34").NumberFormat = "@"
$ws.Range("E34").Value = "  +14.24%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "0.0588"
$ws.Range("E35").NumberFormat = "@"
$ws.Range("E35").Value = "  +4.78%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "4.27"
$ws.Range("E36").NumberFormat = "@"
$ws.Range("E36").Value = "  +5.21%  "
$ws.Range("E37").NumberFormat = "@"
$ws.Range("E37").Value = "  -0.40%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.910"
$ws.Range("E38").NumberFormat = "@"
$ws.Range("E38").Value = "  -2.59%  "
$ws.Range("E39").NumberFormat = "@"
$ws.Range("E39").Value = "  +1.21%  "
$ws.Range("E40").NumberFormat = "@"
$ws.Range("E40").Value = "  +5.94%  "
$ws.Range("E43").NumberFormat = "@"
$ws.Range("E43").Value = "  +3.33%  "
$ws.Range("E44").NumberFormat = "@"
$ws.Range("E44").Value = "  +4.80%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "2.45"
$ws.Range("E45").NumberFormat = "@"
$ws.Range("E45").Value = "  +3.95%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "1.333.73"
$ws.Range("E46").NumberFormat = "@"
$ws.Range("E46").Value = "  +0.05%  "
$ws.Range("E48").NumberFormat = "@"
$ws.Range("E48").Value = "  +4.02%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "2.75"
$ws.Range("E49").NumberFormat = "@"
$ws.Range("E49").Value = "  -0.27%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "45.38"
$ws.Range("E50").NumberFormat = "@"
$ws.Range("E50").Value = "  -5.14%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "0.0754"
$ws.Range("E51").NumberFormat = "@"
$ws.Range("E51").Value = "  +7.16%  "

# Rows 41/42: InjectiveProtocol and Aave swap rank positions,
# each also receiving freshly scraped price/volume data.
$ws.Range("B41").NumberFormat = "@"
$ws.Range("B41").Value = "Aave"
$ws.Range("C41").NumberFormat = "@"
$ws.Range("C41").Value = "https://coinranking.com/coin/ixgUfzmLR+aave-aave"
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "97.14"
$ws.Range("E41").NumberFormat = "@"
$ws.Range("E41").Value = "  +9.28%  "

$ws.Range("B42").NumberFormat = "@"
$ws.Range("B42").Value = "InjectiveProtocol"
$ws.Range("C42").NumberFormat = "@"
$ws.Range("C42").Value = "https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj"
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "17.02"
$ws.Range("E42").NumberFormat = "@"
$ws.Range("E42").Value = "  +6.87%  "
